$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "x"
$ws.Range("A2").Value = 20

$ws.Range("A2").Select()
